$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Outline"): merge the two runs "AWK" + ": when to use it..." into
# a single run.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$awkPara = $tr2.Paragraphs(3, 1)
$awkPara.Text = "TEMP_MERGE"
$tr2.Paragraphs(3, 1).Text = "AWK: when to use it, why people use it so much, why its so important for GMT, and how to use it! (most likely next time)"

# ---------------------------------------------------------------------------
# Slide 3 ("Bash scripting"): consolidate several run splits that the author
# cleaned up, without altering the visible text or formatting.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Paragraph 5: "Makes your life easier: when you are managing data files..."
#   run1 "Makes your "               (unchanged)
#   run2+run3 "life "+"easier:"      -> "life easier:"               (yellow)
#   run4+run5 " "+"when you are..."  -> " when you are managing..."  (plain)
$para5 = $tr3.Paragraphs(5, 1)
$start5 = $para5.Start

$lifeEasier = $tr3.Characters($start5 + 11, 12)
$lifeEasier.Text = "TEMP_MERGE"
$tr3.Characters($start5 + 11, 10).Text = "life easier:"

$whenYou = $tr3.Characters($start5 + 23, 82)
$whenYou.Text = "TEMP_MERGE"
$tr3.Characters($start5 + 23, 10).Text = " when you are managing data files and don" + [char]0x2019 + "t want to waste time doing it " + [char]0x201C + "manually" + [char]0x201D

# Paragraph 8: "When your work requires to integrate different languages: "
#   run2+run3 "integrate different "+"languages:" -> "integrate different languages:" (yellow)
$para8 = $tr3.Paragraphs(8, 1)
$start8 = $para8.Start

$integrate = $tr3.Characters($start8 + 27, 30)
$integrate.Text = "TEMP_MERGE"
$tr3.Characters($start8 + 27, 10).Text = "integrate different languages:"

# Paragraph 9 (lvl 1): "you need to do something ... GMT.." -> single run
$para9 = $tr3.Paragraphs(9, 1)
$para9.Text = "TEMP_MERGE"
$tr3.Paragraphs(9, 1).Text = "you need to do something that first runs a FORTRAN program and then something in C++ or you want to run a program and at the end get a plot made with GMT.."

# ---------------------------------------------------------------------------
# Slide 8 ("Exercise #1"): consolidate run splits in the "Run it" paragraphs.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange

$para5_8 = $tr8.Paragraphs(5, 1)
$para5_8.Text = "TEMP_MERGE"
$tr8.Paragraphs(5, 1).Text = "Run it (just type the name of the file. in some OS" + [char]0x2019 + "s you have to type"

$para6_8 = $tr8.Paragraphs(6, 1)
$para6_8.Text = "TEMP_MERGE"
$tr8.Paragraphs(6, 1).Text = [char]0x201C + "./filename" + [char]0x201D + " to run it)."

# ---------------------------------------------------------------------------
# Remove the five trailing placeholder "AWK" slides (positions 20-24).
# Deleted from the end so earlier indices stay valid.
# ---------------------------------------------------------------------------
for ($i = $p.Slides.Count; $i -ge 20; $i--) {
    $p.Slides.Item($i).Delete()
}
